$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the RANDBETWEEN(10,100) formulas in B2:B15 with a hard-coded
# value of 100 (the formula is gone, only the literal value remains).
# The dependent C2:C15 formulas (=B{row}/1.3) are left untouched - they
# simply recalculate against the new B values.
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 2).Value = 100
}

# Recalculate so the cached formula results for column C reflect the new
# column B values (run_com also auto-recalcs after the script finishes,
# this just makes the intent explicit).
$excel.Calculate() | Out-Null

# Update the active selection/view: active cell B3 with B3:B15 selected.
$ws.Activate() | Out-Null
$ws.Range("B3:B15").Select() | Out-Null
